# Auto-generated Excel COM-interop script
# Applies numeric value updates to specific cells across multiple sheets
# as described by the upstream OOXML diff (scheduled market-price refresh).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 22915.2
$ws.Cells.Item(43, 9).Value = 5500
$ws.Cells.Item(43, 10).Value = 34525.332
$ws.Cells.Item(43, 11).Value = 5500
$ws.Cells.Item(43, 12).Value = 34525.332
$ws.Cells.Item(43, 13).Value = -5431
$ws.Cells.Item(43, 14).Value = -34663.332
$ws.Cells.Item(53, 8).Value = 1773
$ws.Cells.Item(53, 9).Value = 2816.875
$ws.Cells.Item(53, 10).Value = 102.8
$ws.Cells.Item(53, 11).Value = 2816.875
$ws.Cells.Item(53, 12).Value = 102.8
$ws.Cells.Item(53, 13).Value = -2179.875
$ws.Cells.Item(53, 14).Value = -1376.8
$ws.Cells.Item(55, 8).Value = 404.36365
$ws.Cells.Item(55, 9).Value = 399.83334
$ws.Cells.Item(55, 10).Value = 409.8
$ws.Cells.Item(55, 11).Value = 399.83334
$ws.Cells.Item(55, 12).Value = 409.8
$ws.Cells.Item(55, 13).Value = -185.83334
$ws.Cells.Item(55, 14).Value = -837.8
$ws.Cells.Item(74, 8).Value = 3998.6667
$ws.Cells.Item(74, 9).Value = 4246
$ws.Cells.Item(74, 10).Value = 3875
$ws.Cells.Item(74, 11).Value = 4246
$ws.Cells.Item(74, 12).Value = 3875
$ws.Cells.Item(74, 13).Value = -3310
$ws.Cells.Item(74, 14).Value = -5747
$ws.Cells.Item(77, 8).Value = 3998.6667
$ws.Cells.Item(77, 9).Value = 4246
$ws.Cells.Item(77, 10).Value = 3875
$ws.Cells.Item(77, 11).Value = 21230
$ws.Cells.Item(77, 12).Value = 19375
$ws.Cells.Item(77, 13).Value = -16550
$ws.Cells.Item(77, 14).Value = -28735
$ws.Cells.Item(98, 8).Value = 1795.4736
$ws.Cells.Item(98, 9).Value = 1710
$ws.Cells.Item(98, 10).Value = 1872.4
$ws.Cells.Item(98, 11).Value = 1710
$ws.Cells.Item(98, 12).Value = 1872.4
$ws.Cells.Item(98, 13).Value = -212
$ws.Cells.Item(98, 14).Value = -4868.4
$ws.Cells.Item(119, 8).Value = 1000
$ws.Cells.Item(119, 10).Value = 1000
$ws.Cells.Item(119, 12).Value = 3000
$ws.Cells.Item(119, 14).Value = -12676
$ws.Cells.Item(122, 8).Value = 1795.4736
$ws.Cells.Item(122, 9).Value = 1710
$ws.Cells.Item(122, 10).Value = 1872.4
$ws.Cells.Item(122, 11).Value = 5130
$ws.Cells.Item(122, 12).Value = 5617.200000000001
$ws.Cells.Item(122, 13).Value = -2680
$ws.Cells.Item(122, 14).Value = -10517.2
$ws.Cells.Item(139, 8).Value = 41925
$ws.Cells.Item(139, 10).Value = 41925
$ws.Cells.Item(139, 12).Value = 41925
$ws.Cells.Item(139, 14).Value = -52205
$ws.Cells.Item(140, 8).Value = 71142.86
$ws.Cells.Item(140, 10).Value = 71142.86
$ws.Cells.Item(140, 12).Value = 71142.86
$ws.Cells.Item(140, 14).Value = -81502.86
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 8476709
$ws.Cells.Item(74, 9).Value = 11906087
$ws.Cells.Item(74, 10).Value = 4130
$ws.Cells.Item(74, 11).Value = 11906087
$ws.Cells.Item(74, 12).Value = 4130
$ws.Cells.Item(74, 13).Value = -11905213
$ws.Cells.Item(74, 14).Value = -5878
$ws.Cells.Item(77, 8).Value = 8476709
$ws.Cells.Item(77, 9).Value = 11906087
$ws.Cells.Item(77, 10).Value = 4130
$ws.Cells.Item(77, 11).Value = 59530435
$ws.Cells.Item(77, 12).Value = 20650
$ws.Cells.Item(77, 13).Value = -59526067
$ws.Cells.Item(77, 14).Value = -29386
$ws.Cells.Item(132, 8).Value = 6758894.5
$ws.Cells.Item(132, 9).Value = 11906649
$ws.Cells.Item(132, 10).Value = 2465.75
$ws.Cells.Item(132, 11).Value = 35719947
$ws.Cells.Item(132, 12).Value = 7397.25
$ws.Cells.Item(132, 13).Value = -35717417
$ws.Cells.Item(132, 14).Value = -12457.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(24, 8).Value = 1193.5
$ws.Cells.Item(24, 9).Value = 609.6
$ws.Cells.Item(24, 10).Value = 2166.6667
$ws.Cells.Item(24, 11).Value = 609.6
$ws.Cells.Item(24, 12).Value = 2166.6667
$ws.Cells.Item(24, 13).Value = -374.6
$ws.Cells.Item(24, 14).Value = -2636.6667
$ws.Cells.Item(97, 8).Value = 17609.334
$ws.Cells.Item(97, 9).Value = 16414
$ws.Cells.Item(97, 10).Value = 20000
$ws.Cells.Item(97, 11).Value = 16414
$ws.Cells.Item(97, 12).Value = 20000
$ws.Cells.Item(97, 13).Value = -15423
$ws.Cells.Item(97, 14).Value = -21982
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6803873
$ws.Cells.Item(31, 9).Value = 1055
$ws.Cells.Item(31, 11).Value = 1055
$ws.Cells.Item(31, 13).Value = -760
$ws.Cells.Item(34, 8).Value = 6803873
$ws.Cells.Item(34, 9).Value = 1055
$ws.Cells.Item(34, 11).Value = 1055
$ws.Cells.Item(34, 13).Value = -853
$ws.Cells.Item(58, 8).Value = 3688.125
$ws.Cells.Item(58, 9).Value = 1181.3
$ws.Cells.Item(58, 10).Value = 7866.1665
$ws.Cells.Item(58, 11).Value = 1181.3
$ws.Cells.Item(58, 12).Value = 7866.1665
$ws.Cells.Item(58, 13).Value = -978.3
$ws.Cells.Item(58, 14).Value = -8272.166499999999
$ws.Cells.Item(99, 8).Value = 1503.25
$ws.Cells.Item(99, 9).Value = 1482.4
$ws.Cells.Item(99, 10).Value = 1538
$ws.Cells.Item(99, 11).Value = 1482.4
$ws.Cells.Item(99, 12).Value = 1538
$ws.Cells.Item(99, 13).Value = 15.59999999999991
$ws.Cells.Item(99, 14).Value = -4534
$ws.Cells.Item(107, 8).Value = 696.46155
$ws.Cells.Item(107, 9).Value = 645.4
$ws.Cells.Item(107, 10).Value = 866.6667
$ws.Cells.Item(107, 11).Value = 645.4
$ws.Cells.Item(107, 12).Value = 866.6667
$ws.Cells.Item(107, 13).Value = 1274.6
$ws.Cells.Item(107, 14).Value = -4706.6667
$ws.Cells.Item(122, 8).Value = 1569.3448
$ws.Cells.Item(122, 9).Value = 1583.0526
$ws.Cells.Item(122, 10).Value = 1543.3
$ws.Cells.Item(122, 11).Value = 4749.1578
$ws.Cells.Item(122, 12).Value = 4629.9
$ws.Cells.Item(122, 13).Value = -2299.1578
$ws.Cells.Item(122, 14).Value = -9529.9
$ws.Cells.Item(126, 8).Value = 1503.25
$ws.Cells.Item(126, 9).Value = 1482.4
$ws.Cells.Item(126, 10).Value = 1538
$ws.Cells.Item(126, 11).Value = 4447.200000000001
$ws.Cells.Item(126, 12).Value = 4614
$ws.Cells.Item(126, 13).Value = -1977.200000000001
$ws.Cells.Item(126, 14).Value = -9554
$ws.Cells.Item(136, 8).Value = 3688.125
$ws.Cells.Item(136, 9).Value = 1181.3
$ws.Cells.Item(136, 10).Value = 7866.1665
$ws.Cells.Item(136, 11).Value = 3543.9
$ws.Cells.Item(136, 12).Value = 23598.4995
$ws.Cells.Item(136, 13).Value = -993.8999999999996
$ws.Cells.Item(136, 14).Value = -28698.4995
$ws.Cells.Item(140, 8).Value = 26506.666
$ws.Cells.Item(140, 10).Value = 26506.666
$ws.Cells.Item(140, 12).Value = 26506.666
$ws.Cells.Item(140, 14).Value = -36866.666
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 3441.1428
$ws.Cells.Item(5, 9).Value = 517.6
$ws.Cells.Item(5, 11).Value = 1552.8
$ws.Cells.Item(5, 13).Value = -1440.8
$ws.Cells.Item(12, 8).Value = 33.95238
$ws.Cells.Item(12, 9).Value = 12.833333
$ws.Cells.Item(12, 11).Value = 38.499999
$ws.Cells.Item(12, 13).Value = 134.500001
$ws.Cells.Item(80, 8).Value = 4996
$ws.Cells.Item(80, 10).Value = 4996
$ws.Cells.Item(80, 12).Value = 14988
$ws.Cells.Item(80, 14).Value = -16860
$ws.Cells.Item(83, 8).Value = 4996
$ws.Cells.Item(83, 10).Value = 4996
$ws.Cells.Item(83, 12).Value = 44964
$ws.Cells.Item(83, 14).Value = -54324
$ws.Cells.Item(135, 8).Value = 3441.1428
$ws.Cells.Item(135, 9).Value = 517.6
$ws.Cells.Item(135, 11).Value = 4658.400000000001
$ws.Cells.Item(135, 13).Value = -2123.400000000001
$ws.Cells.Item(137, 8).Value = 9808068
$ws.Cells.Item(137, 10).Value = 6787.5
$ws.Cells.Item(137, 12).Value = 20362.5
$ws.Cells.Item(137, 14).Value = -30562.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1092.5625
$ws.Cells.Item(22, 9).Value = 580
$ws.Cells.Item(22, 10).Value = 1325.5454
$ws.Cells.Item(22, 11).Value = 580
$ws.Cells.Item(22, 12).Value = 1325.5454
$ws.Cells.Item(22, 13).Value = -285
$ws.Cells.Item(22, 14).Value = -1915.5454
$ws.Cells.Item(27, 8).Value = 1092.5625
$ws.Cells.Item(27, 9).Value = 580
$ws.Cells.Item(27, 10).Value = 1325.5454
$ws.Cells.Item(27, 11).Value = 580
$ws.Cells.Item(27, 12).Value = 1325.5454
$ws.Cells.Item(27, 13).Value = -473
$ws.Cells.Item(27, 14).Value = -1539.5454
$ws.Cells.Item(139, 8).Value = 41508.547
$ws.Cells.Item(139, 10).Value = 41594.4
$ws.Cells.Item(139, 12).Value = 41594.4
$ws.Cells.Item(139, 14).Value = -51874.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 10000583
$ws.Cells.Item(3, 9).Value = 12000100
$ws.Cells.Item(3, 10).Value = 3000
$ws.Cells.Item(3, 11).Value = 12000100
$ws.Cells.Item(3, 12).Value = 3000
$ws.Cells.Item(3, 13).Value = -11999986
$ws.Cells.Item(3, 14).Value = -3228
$ws.Cells.Item(135, 8).Value = 40833.332
$ws.Cells.Item(135, 10).Value = 40833.332
$ws.Cells.Item(135, 12).Value = 40833.332
$ws.Cells.Item(135, 14).Value = -50973.332
